$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the selection / active cell shown in the sheet view
$ws.Range("I13").Select()

# Row 11 (local-1) updated stats
$ws.Range("G11").Value = 460
$ws.Range("H11").Value = 464
$ws.Range("I11").Value = 492
$ws.Range("J11").Value = 2.4700000000000002
$ws.Range("K11").Value = 464.4
$ws.Range("L11").Value = 30
$ws.Range("M11").Value = 492

# Row 12 (GRASP-1) updated stats
$ws.Range("G12").Value = 474
$ws.Range("H12").Value = 468
$ws.Range("I12").Value = 496
$ws.Range("J12").Value = 3
$ws.Range("K12").Value = 477.2
$ws.Range("L12").Value = 30
$ws.Range("M12").Value = 496

# Row 13 updated stats - G/H updated, I/K/L/M cleared out, J cleared but keeps its style
$ws.Range("G13").Value = 452
$ws.Range("H13").Value = 450
$ws.Range("I13").ClearContents()
$ws.Range("J13").ClearContents()
$ws.Range("K13").ClearContents()
$ws.Range("L13").ClearContents()
$ws.Range("M13").ClearContents()

# Row 14 - newly populated stats
$ws.Range("G14").Value = 2922
$ws.Range("H14").Value = 2934
$ws.Range("I14").Value = 2958
$ws.Range("J14").Value = 3
$ws.Range("K14").Value = 2934
$ws.Range("L14").Value = 30
$ws.Range("M14").Value = 2958
$ws.Range("N14").Value = 3169
